# SOCs for blood and urine samples
# Adds two new Scope-of-Content (SOC) rows to the "SSD" sheet describing
# the Blood Sample and Urine Sample SOCs, mirroring the existing rows.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SSD")
$ws2 = $wb.Worksheets.Item("SOC-NHANES-2017-2018-SUBJECTS")

# --- New row 10: Blood Sample SOC ---
$ws1.Range("B10").Value = "nhanes-kb:SOC-NHANES-2017-2018-BLOOD"
$ws1.Range("D10").Value = "??blood"
$ws1.Range("H10").Value = "Blood Sample"
$ws1.Range("I10").Value = "nhanes-kb:STD-NHANES-2017-2018"
$ws1.Range("J10").Value = "nhanes-kb:SOC-NHANES-2017-2018-SUBJECTS"

# --- New row 11: Urine Sample SOC ---
$ws1.Range("B11").Value = "nhanes-kb:SOC-NHANES-2017-2018-URINE"
$ws1.Range("D11").Value = "??urine"
$ws1.Range("H11").Value = "Urine Sample"
$ws1.Range("I11").Value = "nhanes-kb:STD-NHANES-2017-2018"
$ws1.Range("J11").Value = "nhanes-kb:SOC-NHANES-2017-2018-SUBJECTS"

# Column B widened slightly to fit the new (longer) labels.
$ws1.Columns.Item(2).ColumnWidth = 45.666666666666664

# Selection on the reference subjects sheet moves to A2:A101 without
# making that sheet the active tab.
$ws2.Range("A2:A101").Select()

# Restore "SSD" as the active sheet/tab and set its final selection.
$ws1.Activate()
$ws1.Range("K10").Select()
